# Update "想去人数" (want-to-go count) figures across sheets to match the
# latest generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 5
$ws.Range("F7").Value  = 4549
$ws.Range("F11").Value = 88
$ws.Range("F13").Value = 688
$ws.Range("F15").Value = 970
$ws.Range("F16").Value = 77
$ws.Range("F20").Value = 110
$ws.Range("F21").Value = 93
$ws.Range("F22").Value = 3469
$ws.Range("F23").Value = 5806
$ws.Range("F29").Value = 3346
$ws.Range("F30").Value = 356
$ws.Range("F31").Value = 23
$ws.Range("F32").Value = 2457
$ws.Range("F35").Value = 122
$ws.Range("F36").Value = 208
$ws.Range("F41").Value = 900
$ws.Range("F42").Value = 13
$ws.Range("F45").Value = 44
$ws.Range("F46").Value = 466
$ws.Range("F47").Value = 62

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 95

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 6
$ws.Range("F7").Value  = 4549
$ws.Range("F11").Value = 95
$ws.Range("F12").Value = 88
$ws.Range("F14").Value = 688
$ws.Range("F16").Value = 970
$ws.Range("F17").Value = 77
$ws.Range("F21").Value = 110
$ws.Range("F22").Value = 93
$ws.Range("F23").Value = 3469
$ws.Range("F24").Value = 5806
$ws.Range("F30").Value = 3346
$ws.Range("F31").Value = 356
$ws.Range("F32").Value = 23
$ws.Range("F33").Value = 2457
$ws.Range("F36").Value = 122
$ws.Range("F37").Value = 208
$ws.Range("F42").Value = 901
$ws.Range("F43").Value = 13
$ws.Range("F46").Value = 44
$ws.Range("F47").Value = 466
$ws.Range("F48").Value = 62
